# improved crowding + init by heuristic
# Adds a new "crowding" experiment run (row 5/6 of the benchmark sheet) and
# records its setup description in the Setups lookup table (J/K columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the second run (greedy init + improved crowding) results for tour29
$ws.Range("D5").Value = 2
$ws.Range("E5").Formula = "=AVERAGE(28377, 28907, 29424)"
$ws.Range("F5").Formula = "=AVERAGE(51.79, 62.86, 48.13,38.67)"

# Register the new setup description (setup #2) in the lookup table
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = "popsize=250, offspring=150, k=5, alpha=0.05, prc=0.99, crowding (chance:0.5, hammingdistance)"

# Move the active selection like the author left it
$ws.Range("G20").Select()
